## "modify ddl of ch2 & ch3"
##
## 1) The "date updated" field cached on the Slide Master and on every
##    Slide Layout (a datetimeFigureOut field reading "2021/9/28") moves
##    on three days to "2021/9/30".
## 2) Slide 4's ch2 submission-deadline textbox changes its day-of-month
##    from the 3rd to the 10th ("10 month 3 day" -> "10 month 10 day"),
##    and the textbox (which auto-fits its text) grows to keep its right
##    edge fixed while the left edge moves to accommodate the now two
##    digit day number.
## 3) Slide 7's ch3 submission-deadline textbox changes its day-of-month
##    from the 10th to the 30th ("10 month 10 day" -> "10 month 30 day").
##    Both values are already two digits, so the textbox geometry is
##    unchanged.

$p = $ppt.ActivePresentation

function Set-CachedDate($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "2021/9/28") {
                $sh.TextFrame.TextRange.Text = "2021/9/30"
            }
        }
    }
}

# --- 1) refresh the cached "datetimeFigureOut" text on master + every layout
$master = $p.SlideMaster
Set-CachedDate $master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Set-CachedDate $master.CustomLayouts.Item($li)
}

# --- 2) slide 4 (ch2): deadline day 3 -> 10, textbox grows (right edge pinned)
$slide4 = $p.Slides.Item(4)
for ($i = 1; $i -le $slide4.Shapes.Count; $i++) {
    $sh = $slide4.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "提交截止日期：10月3日24:00") {
            $sh.TextFrame.TextRange.Text = "提交截止日期：10月10日24:00"
            $sh.Left = 648.6991119384766
            $sh.Width = 292.2910690307617
        }
    }
}

# --- 3) slide 7 (ch3): deadline day 10 -> 30 (geometry unchanged)
$slide7 = $p.Slides.Item(7)
for ($i = 1; $i -le $slide7.Shapes.Count; $i++) {
    $sh = $slide7.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "提交截止日期：10月10日24:00") {
            $sh.TextFrame.TextRange.Text = "提交截止日期：10月30日24:00"
        }
    }
}
